# Scheduled-runner price/profit refresh across the Sargatanas_Profits sheets.
# Updates cached market-board price + profit figures (columns H-N) for the
# affected leve rows on each job sheet; no formulas are involved, these are
# static snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4479.4
$ws.Range("I19").Value = 6666.3335
$ws.Range("K19").Value = 6666.3335
$ws.Range("M19").Value = -6491.3335
$ws.Range("H33").Value = 1740
$ws.Range("I33").Value = 2989.2
$ws.Range("J33").Value = 699
$ws.Range("K33").Value = 2989.2
$ws.Range("L33").Value = 699
$ws.Range("M33").Value = -2760.2
$ws.Range("N33").Value = -1157
$ws.Range("H70").Value = 27779772
$ws.Range("I70").Value = 50002156
$ws.Range("J70").Value = 20835276
$ws.Range("K70").Value = 150006468
$ws.Range("L70").Value = 62505828
$ws.Range("M70").Value = -150006198
$ws.Range("N70").Value = -62506368
$ws.Range("H73").Value = 27779772
$ws.Range("I73").Value = 50002156
$ws.Range("J73").Value = 20835276
$ws.Range("K73").Value = 150006468
$ws.Range("L73").Value = 62505828
$ws.Range("M73").Value = -150005532
$ws.Range("N73").Value = -62507700
$ws.Range("H92").Value = 1578.0588
$ws.Range("I92").Value = 762
$ws.Range("J92").Value = 3074.1667
$ws.Range("K92").Value = 762
$ws.Range("L92").Value = 3074.1667
$ws.Range("M92").Value = 486
$ws.Range("N92").Value = -5570.1667
$ws.Range("H96").Value = 1110.375
$ws.Range("I96").Value = 981.5
$ws.Range("J96").Value = 1497
$ws.Range("K96").Value = 2944.5
$ws.Range("L96").Value = 4491
$ws.Range("M96").Value = -1571.5
$ws.Range("N96").Value = -7237
$ws.Range("H103").Value = 1018.82355
$ws.Range("J103").Value = 1129.0714
$ws.Range("L103").Value = 3387.2142
$ws.Range("N103").Value = -4559.2142
$ws.Range("H129").Value = 1107.8334
$ws.Range("I129").Value = 661.75
$ws.Range("K129").Value = 1985.25
$ws.Range("M129").Value = 3014.75
$ws.Range("H135").Value = 400687.03
$ws.Range("I135").Value = 435442.44
$ws.Range("K135").Value = 3918981.96
$ws.Range("M135").Value = -3916446.96
$ws.Range("H137").Value = 3255.923
$ws.Range("I137").Value = 3275.7144
$ws.Range("K137").Value = 9827.143199999999
$ws.Range("M137").Value = -7277.143199999999
$ws.Range("H138").Value = 5111.387
$ws.Range("J138").Value = 11772.637
$ws.Range("L138").Value = 35317.911
$ws.Range("N138").Value = -45597.911
$ws.Range("H141").Value = 14495719
$ws.Range("I141").Value = 16669437
$ws.Range("K141").Value = 50008311
$ws.Range("M141").Value = -50003131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8177.125
$ws.Range("I45").Value = 2991.3333
$ws.Range("J45").Value = 11288.6
$ws.Range("K45").Value = 2991.3333
$ws.Range("L45").Value = 11288.6
$ws.Range("M45").Value = -2614.3333
$ws.Range("N45").Value = -12042.6
$ws.Range("H53").Value = 9900
$ws.Range("J53").Value = 9900
$ws.Range("L53").Value = 9900
$ws.Range("N53").Value = -11264
$ws.Range("J110").Value = 66668190
$ws.Range("L110").Value = 66668190
$ws.Range("N110").Value = -66672280
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 72878.32000000001
$ws.Range("I105").Value = 98522.19
$ws.Range("J105").Value = 4494.6665
$ws.Range("K105").Value = 98522.19
$ws.Range("L105").Value = 4494.6665
$ws.Range("M105").Value = -96775.19
$ws.Range("N105").Value = -7988.6665
$ws.Range("H134").Value = 5259.44
$ws.Range("I134").Value = 1967.963
$ws.Range("K134").Value = 5903.889
$ws.Range("M134").Value = -3368.889
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13340116
$ws.Range("I132").Value = 2056.5
$ws.Range("J132").Value = 22232156
$ws.Range("K132").Value = 6169.5
$ws.Range("L132").Value = 66696468
$ws.Range("M132").Value = -3639.5
$ws.Range("N132").Value = -66701528
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79349.07000000001
$ws.Range("I2").Value = 13085.782
$ws.Range("J2").Value = 333358.34
$ws.Range("K2").Value = 78514.692
$ws.Range("L2").Value = 2000150.04
$ws.Range("M2").Value = -78401.692
$ws.Range("N2").Value = -2000376.04
$ws.Range("H94").Value = 4300
$ws.Range("I94").Value = 3250
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 9750
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = -9074
$ws.Range("N94").Value = -16352
$ws.Range("H121").Value = 1164.4286
$ws.Range("J121").Value = 1997.75
$ws.Range("L121").Value = 5993.25
$ws.Range("N121").Value = -8613.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 89995
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 89995
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 89995
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -90513
$ws.Range("H70").Value = 7179.6
$ws.Range("I70").Value = 7115.2856
$ws.Range("J70").Value = 7329.6665
$ws.Range("K70").Value = 7115.2856
$ws.Range("L70").Value = 7329.6665
$ws.Range("M70").Value = -6845.2856
$ws.Range("N70").Value = -7869.6665
$ws.Range("H73").Value = 7179.6
$ws.Range("I73").Value = 7115.2856
$ws.Range("J73").Value = 7329.6665
$ws.Range("K73").Value = 7115.2856
$ws.Range("L73").Value = 7329.6665
$ws.Range("M73").Value = -6179.2856
$ws.Range("N73").Value = -9201.666499999999
$ws.Range("H80").Value = 4353.727
$ws.Range("I80").Value = 2799.6
$ws.Range("K80").Value = 2799.6
$ws.Range("M80").Value = -1801.6
$ws.Range("H83").Value = 4353.727
$ws.Range("I83").Value = 2799.6
$ws.Range("K83").Value = 13998
$ws.Range("M83").Value = -9006
$ws.Range("H102").Value = 3203.5454
$ws.Range("I102").Value = 3273.889
$ws.Range("K102").Value = 3273.889
$ws.Range("M102").Value = -1651.889
$ws.Range("H113").Value = 2958.0667
$ws.Range("I113").Value = 3043.6667
$ws.Range("K113").Value = 3043.6667
$ws.Range("M113").Value = -873.6667000000002
$ws.Range("H126").Value = 7949.5
$ws.Range("I126").Value = 7949.5
$ws.Range("K126").Value = 23848.5
$ws.Range("M126").Value = -21378.5
$ws.Range("H132").Value = 4595
$ws.Range("I132").Value = 1596.5625
$ws.Range("K132").Value = 4789.6875
$ws.Range("M132").Value = -2259.6875
